$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24, shifting existing rows 24-61 down to 25-62
$ws.Range("A24").EntireRow.Insert()

# Populate the new row 24 with the new entry
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = '2012'
$ws.Range("A24").Style = "Normal"

$ws.Range("B24").Value = '**** <br> [Antibiotic Resistance in Bacterial Pathogens from Retail Raw Meats and Food-Producing Animals in Japan](https://www.sciencedirect.com/science/article/pii/S0362028X23039777?via%3Dihub) <br> (Journal of Food Protection, Volume 75, Issue 10, 1 October 2012, Pages 1774-1782)'
$ws.Range("C24").Value = '未登録'
